# Update "想去人数" (want-to-go count) values in column F across the four sheets
# to match the refreshed data snapshot ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$changes = @{
    "展览" = @{
        2 = 16
        3 = 975
        6 = 1090
        7 = 875
        8 = 265
        9 = 59
        11 = 860
        12 = 301
        15 = 1353
        17 = 1275
        18 = 1226
        19 = 2906
        21 = 730
        25 = 1044
        26 = 365
        27 = 3215
        28 = 629
        29 = 540
        30 = 1430
    }
    "演出" = @{
        5 = 61
        8 = 30
        12 = 7
    }
    "本地生活" = @{
        2 = 758
    }
    "全部类型" = @{
        2 = 16
        3 = 758
        6 = 975
        10 = 1090
        11 = 875
        12 = 265
        13 = 61
        14 = 59
        18 = 30
        21 = 860
        22 = 301
        25 = 1353
        27 = 1275
        28 = 1226
        29 = 2906
        31 = 730
        37 = 1044
        38 = 365
        39 = 3215
        40 = 629
        41 = 540
        42 = 1430
        43 = 7
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowChanges = $changes[$sheetName]
    foreach ($row in $rowChanges.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowChanges[$row]
    }
}
